$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy original column D (which will become the source format for the two new
# columns) and insert-copy it twice before D, shifting everything right.
$ws.Columns("D").Copy() | Out-Null
$ws.Columns("D").Insert() | Out-Null
$ws.Columns("D").Copy() | Out-Null
$ws.Columns("D").Insert() | Out-Null
